$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.183.88'
$ws.Range("E2").Value = '  +1.11%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.216.65'
$ws.Range("E3").Value = '  -0.55%  '

# Row 4
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.74'
$ws.Range("E5").Value = '  +0.83%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.72'
$ws.Range("E7").Value = '  -0.96%  '

# Row 8
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.401'
$ws.Range("E9").Value = '  -0.26%  '

# Row 10
$ws.Range("E10").Value = '  +2.33%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.103'
$ws.Range("E11").Value = '  -0.15%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.544.16'
$ws.Range("E12").Value = '  -0.61%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.43'
$ws.Range("E13").Value = '  -1.43%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.11'
$ws.Range("E14").Value = '  +2.58%  '

# Row 15
$ws.Range("E15").Value = '  +0.19%  '

# Row 16
$ws.Range("E16").Value = '  +0.27%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.213.24'
$ws.Range("E17").Value = '  -0.62%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.111.55'
$ws.Range("E18").Value = '  +1.15%  '

# Row 19
$ws.Range("E19").Value = '  +4.29%  '

# Row 20
$ws.Range("E20").Value = '  +2.51%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.05'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '243.93'
$ws.Range("E22").Value = '  -1.81%  '

# Row 23
$ws.Range("E23").Value = '  -0.19%  '

# Row 24
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.39'
$ws.Range("E24").Value = '  +3.36%  '

# Row 25
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.40'
$ws.Range("E25").Value = '  +1.44%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.59'
$ws.Range("E26").Value = '  +0.24%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.26'
$ws.Range("E27").Value = '  +0.86%  '

# Row 28
$ws.Range("E28").Value = '  +0.28%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.32'
$ws.Range("E29").Value = '  +2.01%  '

# Row 30
$ws.Range("E30").Value = '  +2.92%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.65'
$ws.Range("E31").Value = '  +2.38%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.120'
$ws.Range("E32").Value = '  -1.34%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.95'
$ws.Range("E33").Value = '  -1.97%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.60'
$ws.Range("E34").Value = '  -0.92%  '

# Row 35
$ws.Range("E35").Value = '  +3.82%  '

# Row 36
$ws.Range("B36").Value = 'THORChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.31'
$ws.Range("E36").Value = '  -3.88%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.55'
$ws.Range("E37").Value = '  -3.93%  '

# Row 38
$ws.Range("E38").Value = '  -1.48%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0250'
$ws.Range("E39").Value = '  +5.75%  '

# Row 40
$ws.Range("E40").Value = '  -0.11%  '

# Row 41
$ws.Range("E41").Value = '  -3.31%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.52'
$ws.Range("E42").Value = '  -2.39%  '

# Row 43
$ws.Range("E43").Value = '  -2.68%  '

# Row 44
$ws.Range("E44").Value = '  +1.14%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '96.85'
$ws.Range("E45").Value = '  -2.36%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.33'

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.452.30'
$ws.Range("E47").Value = '  -1.16%  '

# Row 48
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.07'
$ws.Range("E48").Value = '  -0.88%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '15.97'
$ws.Range("E49").Value = '  -3.00%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.72'
$ws.Range("E50").Value = '  -2.41%  '

# Row 51
$ws.Range("E51").Value = '  +1.38%  '
